$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Middle Name" column (column C) entirely, shifting all columns
# from D onward one position to the left. This matches Excel's native
# "Delete Column" behavior: cell values/styles shift left, column widths
# shift left, conditional formatting / hyperlinks / drawings auto-adjust.
$ws.Columns.Item(3).Delete()
